$wb = $excel.ActiveWorkbook

# --- Details sheet: move the cursor/selection from K16 to N1 ---
$wsDetails = $wb.Worksheets.Item("Details")
$wsDetails.Activate()
$wsDetails.Range("N1").Select() | Out-Null

# --- Etape_1 sheet: move the selection from C18:D18 to A23:D23 ---
$wsEtape1 = $wb.Worksheets.Item("Etape_1")
$wsEtape1.Activate()
$wsEtape1.Range("A23:D23").Select() | Out-Null

# --- Etape_2 sheet: insert a new waypoint row (row 33) ---
$wsEtape2 = $wb.Worksheets.Item("Etape_2")
$wsEtape2.Activate()

$wsEtape2.Rows.Item(33).Insert()

$wsEtape2.Range("A33").Value = 129.5
$wsEtape2.Range("B33").Value = "Left"
$wsEtape2.Range("C33").Value = "Entrée sur le circuit d'arrivée <br/>14e Avenue E"
$wsEtape2.Range("D33").Value = "Circuit finish entrance<br/>14e Avenue E"

# the waypoint that used to be at 129.5 (now pushed down to row 34) is
# renumbered to 129.6 to keep the km sequence increasing
$wsEtape2.Range("A34").Value = 129.6

$wsEtape2.Range("A33:D33").Select() | Out-Null

# --- Etape_6 sheet: becomes the active tab, with row 7 selected ---
$wsEtape6 = $wb.Worksheets.Item("Etape_6")
$wsEtape6.Activate()
$wsEtape6.Rows.Item(7).Select() | Out-Null
